$wb = $excel.ActiveWorkbook

# --- Sheet "Y" (sheet1.xml) ---
$wsY = $wb.Worksheets.Item("Y")

# Add new data in row 2: B2 = "Maize (home consumed)", C2 = 100000000000000
$wsY.Range("B2").Style = "Normal"
$wsY.Range("B2").Value = "Maize (home consumed)"
$wsY.Range("C2").Value = 100000000000000

# Clear the values previously in A3:A6 (keep formatting/style)
$wsY.Range("A3").ClearContents()
$wsY.Range("A4").ClearContents()
$wsY.Range("A5").ClearContents()
$wsY.Range("A6").ClearContents()

# Update selection on sheet Y to B3 and make it the active/tab-selected sheet
$wsY.Activate()
$wsY.Range("B3").Select()

# --- Sheet "A" (sheet2.xml) ---
$wsA = $wb.Worksheets.Item("A")
$wsA.Activate()
$wsA.Range("F2").Select()

# --- Sheet "VA" (sheet3.xml) ---
$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Activate()
$wsVA.Range("D10").Select()

# Restore the originally-active sheet to "Y" (so tabSelected/activeTab ends up there)
$wsY.Activate()
